$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old rows 2 and 3 (their data is superseded by shifting subsequent rows up)
$ws.Rows("2:3").Delete()

# Add the new row of data at the bottom (row 21 after shifting)
$ws.Range("A21").Value = 0.1471566200256338
$ws.Range("B21").Value = 1.182808732986451
$ws.Range("C21").Value = 0.03668105900287391
